$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2 through 14
# from serial date 45175 (2023-09-06) to 45183 (2023-09-14)
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
